$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Remis"
$ws.Range("C6").Value = "Remis"
$ws.Range("C10").Value = "Remis"
$ws.Range("C11").Value = "Cracovia"
$ws.Range("C12").Value = "Remis"
$ws.Range("C13").Value = "Remis"
$ws.Range("C14").Value = "Lechia Gdańsk"
$ws.Range("C15").Value = "Legia Warszawa"
$ws.Range("C17").Value = "Remis"
$ws.Range("C21").Value = "Zagłębie Lubin"
$ws.Range("C22").Value = "Śląsk Wrocław"
$ws.Range("C24").Value = "Warta Poznań"
$ws.Range("C26").Value = "Remis"
$ws.Range("C27").Value = "Stal Mielec"
$ws.Range("C29").Value = "Remis"
$ws.Range("C33").Value = "Górnik Zabrze"
$ws.Range("C34").Value = "Remis"
$ws.Range("C35").Value = "Warta Poznań"
$ws.Range("C36").Value = "Miedź Legnica"
$ws.Range("C39").Value = "Remis"
$ws.Range("C40").Value = "Lech Poznań"
$ws.Range("C41").Value = "Remis"
$ws.Range("C47").Value = "Piast Gliwice"
$ws.Range("C50").Value = "Lechia Gdańsk"
$ws.Range("C54").Value = "Wisła Płock"
$ws.Range("C59").Value = "Górnik Zabrze"
$ws.Range("C63").Value = "Wisła Płock"
$ws.Range("C65").Value = "Raków Częstochowa"
$ws.Range("C66").Value = "Piast Gliwice"
$ws.Range("C67").Value = "Zagłębie Lubin"
$ws.Range("C68").Value = "Pogoń Szczecin"
$ws.Range("C70").Value = "Remis"
$ws.Range("C75").Value = "Miedź Legnica"
$ws.Range("C76").Value = "Piast Gliwice"
$ws.Range("C79").Value = "Raków Częstochowa"
$ws.Range("C81").Value = "Widzew Łódź"
$ws.Range("C82").Value = "Remis"
$ws.Range("C84").Value = "Śląsk Wrocław"
$ws.Range("C85").Value = "Wisła Płock"
$ws.Range("C87").Value = "Jagielonia Białystok"
$ws.Range("C88").Value = "Legia Warszawa"
$ws.Range("C90").Value = "Remis"
$ws.Range("C93").Value = "Lech Poznań"
$ws.Range("C94").Value = "Remis"
$ws.Range("C95").Value = "Górnik Zabrze"
$ws.Range("C97").Value = "Radomiak Radom"
$ws.Range("C98").Value = "Widzew Łódź"
$ws.Range("C100").Value = "Śląsk Wrocław"
$ws.Range("C104").Value = "Remis"
$ws.Range("C106").Value = "Legia Warszawa"
$ws.Range("C109").Value = "Śląsk Wrocław"
$ws.Range("C110").Value = "Remis"
$ws.Range("C111").Value = "Raków Częstochowa"
$ws.Range("C113").Value = "Lech Poznań"
$ws.Range("C115").Value = "Remis"
$ws.Range("C117").Value = "Zagłębie Lubin"
$ws.Range("C120").Value = "Piast Gliwice"
$ws.Range("C122").Value = "Remis"
$ws.Range("C124").Value = "Wisła Płock"
$ws.Range("C126").Value = "Remis"
$ws.Range("C128").Value = "Remis"
$ws.Range("C131").Value = "Raków Częstochowa"
$ws.Range("C132").Value = "Remis"
$ws.Range("C135").Value = "Radomiak Radom"
$ws.Range("C136").Value = "Wisła Płock"
$ws.Range("C137").Value = "Remis"
$ws.Range("C144").Value = "Zagłębie Lubin"
$ws.Range("C145").Value = "Widzew Łódź"
$ws.Range("C146").Value = "Remis"
$ws.Range("C151").Value = "Remis"
$ws.Range("C152").Value = "Remis"
$ws.Range("C156").Value = "Piast Gliwice"
$ws.Range("C161").Value = "Raków Częstochowa"
$ws.Range("C163").Value = "Zagłębie Lubin"
$ws.Range("C164").Value = "Legia Warszawa"
$ws.Range("C165").Value = "Remis"
$ws.Range("C167").Value = "Lechia Gdańsk"
$ws.Range("C169").Value = "Radomiak Radom"
$ws.Range("C172").Value = "Remis"
$ws.Range("C173").Value = "Remis"
$ws.Range("C174").Value = "Zagłębie Lubin"
$ws.Range("C177").Value = "Remis"
$ws.Range("C180").Value = "Wisła Płock"
$ws.Range("C183").Value = "Piast Gliwice"
$ws.Range("C186").Value = "Wisła Płock"
$ws.Range("C188").Value = "Radomiak Radom"
$ws.Range("C189").Value = "Raków Częstochowa"
$ws.Range("C191").Value = "Cracovia"
$ws.Range("C192").Value = "Jagielonia Białystok"
$ws.Range("C196").Value = "Górnik Zabrze"
$ws.Range("C200").Value = "Śląsk Wrocław"
$ws.Range("C203").Value = "Remis"
$ws.Range("C205").Value = "Pogoń Szczecin"
$ws.Range("C207").Value = "Remis"
$ws.Range("C208").Value = "Remis"
$ws.Range("C211").Value = "Remis"
$ws.Range("C212").Value = "Remis"
$ws.Range("C218").Value = "Remis"
$ws.Range("C219").Value = "Piast Gliwice"
$ws.Range("C220").Value = "Remis"
$ws.Range("C222").Value = "Radomiak Radom"
$ws.Range("C223").Value = "Cracovia"
$ws.Range("C226").Value = "Remis"
$ws.Range("C227").Value = "Remis"
$ws.Range("C228").Value = "Piast Gliwice"
$ws.Range("C229").Value = "Zagłębie Lubin"
$ws.Range("C230").Value = "Miedź Legnica"
$ws.Range("C235").Value = "Radomiak Radom"
$ws.Range("C236").Value = "Jagielonia Białystok"
$ws.Range("C237").Value = "Lech Poznań"
$ws.Range("C240").Value = "Remis"
$ws.Range("C241").Value = "Remis"
$ws.Range("C243").Value = "Remis"
$ws.Range("C244").Value = "Remis"
$ws.Range("C245").Value = "Radomiak Radom"
$ws.Range("C246").Value = "Wisła Płock"
$ws.Range("C248").Value = "Remis"
$ws.Range("C249").Value = "Remis"
$ws.Range("C253").Value = "Warta Poznań"
$ws.Range("C254").Value = "Zagłębie Lubin"
$ws.Range("C257").Value = "Remis"
$ws.Range("C258").Value = "Stal Mielec"
$ws.Range("C260").Value = "Legia Warszawa"
$ws.Range("C261").Value = "Remis"
$ws.Range("C262").Value = "Remis"
$ws.Range("C263").Value = "Cracovia"
$ws.Range("C264").Value = "Remis"
$ws.Range("C266").Value = "Widzew Łódź"
$ws.Range("C268").Value = "Legia Warszawa"
$ws.Range("C271").Value = "Śląsk Wrocław"
$ws.Range("C273").Value = "Remis"
$ws.Range("C275").Value = "Lechia Gdańsk"
$ws.Range("C276").Value = "Widzew Łódź"
$ws.Range("C279").Value = "Piast Gliwice"
$ws.Range("C280").Value = "Stal Mielec"
$ws.Range("C282").Value = "Piast Gliwice"
$ws.Range("C283").Value = "Legia Warszawa"
$ws.Range("C287").Value = "Remis"
$ws.Range("C288").Value = "Górnik Zabrze"
$ws.Range("C289").Value = "Śląsk Wrocław"
$ws.Range("C290").Value = "Remis"
$ws.Range("C291").Value = "Zagłębie Lubin"
$ws.Range("C296").Value = "Piast Gliwice"
$ws.Range("C297").Value = "Remis"
$ws.Range("C299").Value = "Cracovia"
$ws.Range("C300").Value = "Lechia Gdańsk"
$ws.Range("C303").Value = "Remis"
$ws.Range("C307").Value = "Widzew Łódź"
